$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels for the team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) used across row 1
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the team record (Wins/Losses/Ties) for every player row
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 84  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 78  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
